$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the "Epoch N finished" wording -> "Epoch N" in the existing train-log cell (B2)
$ws.Range("B2").Value = "Epoch 1 | Train Loss: 0.452 | Valid Loss: 0.271 | Combined metric: 0.574 | Dice: 0.179 (LB 0.505, SB 0.236, S 0.236) | Hausdorff: 0.162 (LB 0.153, SB 0.095, S 0.239)`n" + `
"Epoch 2 | Train Loss: 0.167 | Valid Loss: 0.280 | Combined metric: 0.540 | Dice: 0.205 (LB 0.499, SB 0.285, S 0.147) | Hausdorff: 0.237 (LB 0.194, SB 0.160, S 0.357)`n" + `
"Epoch 3 | Train Loss: 0.114 | Valid Loss: 0.181 | Combined metric: 0.616 | Dice: 0.241 (LB 0.622, SB 0.207, S 0.561) | Hausdorff: 0.135 (LB 0.174, SB 0.099, S 0.131)`n" + `
"Epoch 4 | Train Loss: 0.091 | Valid Loss: 0.167 | Combined metric: 0.636 | Dice: 0.255 (LB 0.648, SB 0.265, S 0.578) | Hausdorff: 0.110 (LB 0.126, SB 0.132, S 0.072)`n" + `
"Epoch 5 | Train Loss: 0.079 | Valid Loss: 0.170 | Combined metric: 0.646 | Dice: 0.253 (LB 0.649, SB 0.237, S 0.605) | Hausdorff: 0.091 (LB 0.083, SB 0.089, S 0.101)"

# 2) Add the new row (row 3) describing the Unet Efficientnet-b1 model trained on the full train data
$ws.Range("A3").Value = "Unet Efficientnet-b1 `nTrained on full train data"

$ws.Range("B3").Value = "Epoch 1 | Train Loss: 0.203 | Valid Loss: 0.123 | Combined metric: 0.654 | Dice: 0.292 (LB 0.675, SB 0.506, S 0.660) | Hausdorff: 0.104 (LB 0.093, SB 0.164, S 0.054)`n" + `
"Epoch 2 | Train Loss: 0.108 | Valid Loss: 0.114 | Combined metric: 0.671 | Dice: 0.304 (LB 0.688, SB 0.566, S 0.693) | Hausdorff: 0.084 (LB 0.108, SB 0.073, S 0.071)`n" + `
"Epoch 3 | Train Loss: 0.096 | Valid Loss: 0.110 | Combined metric: 0.653 | Dice: 0.304 (LB 0.612, SB 0.565, S 0.694) | Hausdorff: 0.115 (LB 0.214, SB 0.076, S 0.054)`n" + `
"Epoch 4 | Train Loss: 0.087 | Valid Loss: 0.115 | Combined metric: 0.668 | Dice: 0.300 (LB 0.687, SB 0.524, S 0.670) | Hausdorff: 0.087 (LB 0.109, SB 0.096, S 0.057)`n" + `
"Epoch 5 | Train Loss: 0.081 | Valid Loss: 0.107 | Combined metric: 0.689 | Dice: 0.311 (LB 0.703, SB 0.581, S 0.728) | Hausdorff: 0.060 (LB 0.068, SB 0.065, S 0.047)"

$ws.Range("C3").Value = 0.689
$ws.Range("D3").Value = 0.82808
$ws.Range("E3").Value = 0.81848

# Match the wrapped / multi-line formatting used by row 2, and the row height it produces
$ws.Range("A3:B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 57.45

# Update the active selection to reflect where the user ended up after editing
$ws.Range("B11").Select() | Out-Null
